$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.868.52"
$ws.Range("E2").Value = "  -4.53%  "

$ws.Range("D3").Value = "1.736.02"
$ws.Range("E3").Value = "  -4.51%  "

$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.15"
$ws.Range("E5").Value = "  -3.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5748"
$ws.Range("E6").Value = "  -3.02%  "

$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2730"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.01"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06613"
$ws.Range("E10").Value = "  -3.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07541"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").Value = "1.735.43"
$ws.Range("E12").Value = "  -4.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.695"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5993"
$ws.Range("E14").Value = "  -3.30%  "

$ws.Range("D15").Value = "1.973.12"
$ws.Range("E15").Value = "  -4.51%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.34"
$ws.Range("E16").Value = "  -1.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008656"
$ws.Range("E17").Value = "  -9.81%  "

$ws.Range("D18").Value = "27.856.53"
$ws.Range("E18").Value = "  -3.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.304"
$ws.Range("E19").Value = "  -2.83%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "204.71"
$ws.Range("E21").Value = "  -2.09%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.25"
$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.596"
$ws.Range("E23").Value = "  -2.40%  "

$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.89"
$ws.Range("E25").Value = "  -3.21%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.033"
$ws.Range("E26").Value = "  +2.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1227"
$ws.Range("E27").Value = "  -3.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.16"
$ws.Range("E28").Value = "  -0.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.381"
$ws.Range("E29").Value = "  -2.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06166"
$ws.Range("E30").Value = "  -4.10%  "

$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.734"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.721"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.672"
$ws.Range("E34").Value = "  -1.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.032"
$ws.Range("E35").Value = "  -3.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6424"
$ws.Range("E36").Value = "  +1.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.425"
$ws.Range("E37").Value = "  -4.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.706"
$ws.Range("E38").Value = "  -2.07%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01667"
$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("D40").Value = "1.119.36"
$ws.Range("E40").Value = "  -1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.161"
$ws.Range("E41").Value = "  -5.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8723"
$ws.Range("E42").Value = "  -0.75%  "

$ws.Range("E43").Value = "  -0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.88"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("D45").Value = "1.884.42"
$ws.Range("E45").Value = "  -4.82%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.23"
$ws.Range("E46").Value = "  -3.08%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.566"
$ws.Range("E47").Value = "  -1.64%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000107"
$ws.Range("E48").Value = "  -7.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.230"
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05373"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4413"
$ws.Range("E51").Value = "  -2.88%  "
